# Updated cryptos list (price/volume refresh) - GitHub Actions data pull.
# Price (col D) and Volume(1h) (col E) are plain text cells, e.g. "27.268.22"
# and "  +0.17%  ". For numeric-looking price strings we force text entry
# (NumberFormat "@" then ClearFormats) so Excel doesn't coerce them into
# floating point numbers / scientific notation and so the cell keeps its
# original (default) style index.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.268.22'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '1.908.66'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.67'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5323'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.32%  '
$ws.Range('E8').Value = '  +1.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07286'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.07'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9021'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08198'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '95.80'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.98%  '
$ws.Range('E14').Value = '  +1.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.001'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008656'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.81'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.77%  '
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').Value = '27.309.21'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').Value = '1.182.78'
$ws.Range('E20').Value = '  -37.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.049'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.78%  '
$ws.Range('E22').Value = '  +1.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.516'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '150.01'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.290'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.25'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.741'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '116.98'
$ws.Range('D28').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.824'
$ws.Range('D29').ClearFormats()
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.814'
$ws.Range('D30').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09293'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.8383'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +5.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05063'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.226'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.999'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.357'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.696'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5755'
$ws.Range('D38').ClearFormats()
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('E40').Value = '  -0.28%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.295'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +3.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.562'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '117.31'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('E44').Value = '  +0.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4927'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.45%  '
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.17'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.637'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '38.62'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06145'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.56'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.46%  '
